# Generate Report for Handoff
# Adds two new tracked files (cb6d05a5-... and fb0225ff-...) to the
# localization status report, on all three sheets (Overview, zh-cn, de-de),
# pushing the existing ".localization-config" row down by two rows, and
# wires up the matching hyperlinks / shared values.

$wb = $excel.ActiveWorkbook

$HYPER_COLOR = 15570276   # BGR for #6495ED (cornflower blue), matches workbook's HyperLink style

function Style-HyperlinkRange($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $HYPER_COLOR
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop all existing hyperlinks so we can rebuild them (and their relationship
# ids) in the exact final order.
$ws1.Range("A1").Hyperlinks.Delete()

# Make room for the two new rows; the ".localization-config" row (currently
# row 4) shifts down to row 6, carrying its formatting with it.
$ws1.Range("A4:C5").Insert(-4121)

$ws1.Range("A4").Value = "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = "fb0225ff-7798-49a0-b895-f15daf28a9d2.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

# Rebuild hyperlinks in final order (table rel is rId1, so hyperlinks start at rId2).
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/2a0b4a74-ad8b-4746-a694-09018d583dab.md", "", "", "2a0b4a74-ad8b-4746-a694-09018d583dab.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/f33bd561-1cc5-4634-b6d9-074840cd5c87.md", "", "", "f33bd561-1cc5-4634-b6d9-074840cd5c87.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md", "", "", "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/fb0225ff-7798-49a0-b895-f15daf28a9d2.md", "", "", "fb0225ff-7798-49a0-b895-f15daf28a9d2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/.localization-config", "", "", ".localization-config") | Out-Null

Style-HyperlinkRange $ws1.Range("A2:A6")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Hyperlinks.Delete()

$ws2.Range("A4:I5").Insert(-4121)

$ws2.Range("A4").Value = "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.71fa8ddafa44f2a1c11115a532e0472a3cf1ab7b.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-02-26 05:17:34"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Include"

$ws2.Range("A5").Value = "fb0225ff-7798-49a0-b895-f15daf28a9d2.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "fb0225ff-7798-49a0-b895-f15daf28a9d2.e76d4068d00d14b40cad0fdfb56b2a114c63659b.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-02-26 05:17:34"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/2a0b4a74-ad8b-4746-a694-09018d583dab.md", "", "", "2a0b4a74-ad8b-4746-a694-09018d583dab.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/523a85f599840bdc7ff9e34e707ee9bb67987414/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/2a0b4a74-ad8b-4746-a694-09018d583dab.61e9710e3d978abb3f4d242fbc255632186d20fa.zh-cn.xlf", "", "", "2a0b4a74-ad8b-4746-a694-09018d583dab.61e9710e3d978abb3f4d242fbc255632186d20fa.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/f33bd561-1cc5-4634-b6d9-074840cd5c87.md", "", "", "f33bd561-1cc5-4634-b6d9-074840cd5c87.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/523a85f599840bdc7ff9e34e707ee9bb67987414/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/f33bd561-1cc5-4634-b6d9-074840cd5c87.6d954f09a290621ce640595f7db086e52aa7e565.zh-cn.xlf", "", "", "f33bd561-1cc5-4634-b6d9-074840cd5c87.6d954f09a290621ce640595f7db086e52aa7e565.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md", "", "", "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/523a85f599840bdc7ff9e34e707ee9bb67987414/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.71fa8ddafa44f2a1c11115a532e0472a3cf1ab7b.zh-cn.xlf", "", "", "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.71fa8ddafa44f2a1c11115a532e0472a3cf1ab7b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/fb0225ff-7798-49a0-b895-f15daf28a9d2.md", "", "", "fb0225ff-7798-49a0-b895-f15daf28a9d2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/523a85f599840bdc7ff9e34e707ee9bb67987414/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/fb0225ff-7798-49a0-b895-f15daf28a9d2.e76d4068d00d14b40cad0fdfb56b2a114c63659b.zh-cn.xlf", "", "", "fb0225ff-7798-49a0-b895-f15daf28a9d2.e76d4068d00d14b40cad0fdfb56b2a114c63659b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/.localization-config", "", "", ".localization-config") | Out-Null

Style-HyperlinkRange $ws2.Range("A2:A6")
Style-HyperlinkRange $ws2.Range("C2:C5")

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1").Hyperlinks.Delete()

$ws3.Range("A4:I5").Insert(-4121)

$ws3.Range("A4").Value = "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.71fa8ddafa44f2a1c11115a532e0472a3cf1ab7b.de-de.xlf"
$ws3.Range("D4").Value = "2016-02-26 05:17:46"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Include"

$ws3.Range("A5").Value = "fb0225ff-7798-49a0-b895-f15daf28a9d2.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "fb0225ff-7798-49a0-b895-f15daf28a9d2.e76d4068d00d14b40cad0fdfb56b2a114c63659b.de-de.xlf"
$ws3.Range("D5").Value = "2016-02-26 05:17:46"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/2a0b4a74-ad8b-4746-a694-09018d583dab.md", "", "", "2a0b4a74-ad8b-4746-a694-09018d583dab.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b925557caf58b788f2a3c6c2a8dd693ee92f7cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/2a0b4a74-ad8b-4746-a694-09018d583dab.61e9710e3d978abb3f4d242fbc255632186d20fa.de-de.xlf", "", "", "2a0b4a74-ad8b-4746-a694-09018d583dab.61e9710e3d978abb3f4d242fbc255632186d20fa.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/f33bd561-1cc5-4634-b6d9-074840cd5c87.md", "", "", "f33bd561-1cc5-4634-b6d9-074840cd5c87.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b925557caf58b788f2a3c6c2a8dd693ee92f7cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/f33bd561-1cc5-4634-b6d9-074840cd5c87.6d954f09a290621ce640595f7db086e52aa7e565.de-de.xlf", "", "", "f33bd561-1cc5-4634-b6d9-074840cd5c87.6d954f09a290621ce640595f7db086e52aa7e565.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md", "", "", "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b925557caf58b788f2a3c6c2a8dd693ee92f7cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.71fa8ddafa44f2a1c11115a532e0472a3cf1ab7b.de-de.xlf", "", "", "cb6d05a5-2ce3-4d61-ae4e-e8b8add5f4f0.71fa8ddafa44f2a1c11115a532e0472a3cf1ab7b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/e2e/fb0225ff-7798-49a0-b895-f15daf28a9d2.md", "", "", "fb0225ff-7798-49a0-b895-f15daf28a9d2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b925557caf58b788f2a3c6c2a8dd693ee92f7cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/fb0225ff-7798-49a0-b895-f15daf28a9d2.e76d4068d00d14b40cad0fdfb56b2a114c63659b.de-de.xlf", "", "", "fb0225ff-7798-49a0-b895-f15daf28a9d2.e76d4068d00d14b40cad0fdfb56b2a114c63659b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/5ee738aa39cf353527554549dd4a5ed87efa8009/.localization-config", "", "", ".localization-config") | Out-Null

Style-HyperlinkRange $ws3.Range("A2:A6")
Style-HyperlinkRange $ws3.Range("C2:C5")

Write-Host "Report regenerated for handoff: added cb6d05a5-... and fb0225ff-... rows to all sheets."
